$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their original text representation (avoid Excel
# auto-converting numeric-looking strings into actual numbers).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.041.67'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.57%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.820.52'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.89%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.98%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.95'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.66%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.79%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4219'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -2.18%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3661'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.41%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07188'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.28%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8392'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.57%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.71'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -4.41%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.827.67'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.48%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.653'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.31%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07076'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.90%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.277'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.14%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.07'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.09%  '

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.11%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008742'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.01%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.86%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.88'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.13%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.117.49'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.32%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.125'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.85%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.84'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.27%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.048.57'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.41%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.982'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.36%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.08'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.16%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.260'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.04%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.26'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.24%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.268'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.27%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.11'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.75%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08712'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.67%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.175'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -4.96%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7357'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -5.74%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.914'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.41%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.411'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.34%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.000'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.13%  '

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.95%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01947'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.06%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05237'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.31%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.328'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.03%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.870'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.17%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1685'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.64%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5029'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.30%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.543'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.42%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.52'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.56%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '106.08'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.89%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4701'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.02%  '

# Row 48
$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.0000'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.10%  '

# Row 49
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06336'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.27%  '

# Row 50
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.894'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.30%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.645'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.00%  '
